$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the "Förändrad" (C) column date for all data rows (row 2 .. row 398) from 45188 -> 45189
for ($r = 2; $r -le 398; $r++) {
    $ws.Cells.Item($r, 3).Value = 45189
}

# --- 2. Row 4 / Row 5 content change.
# A new logging notice "A 59219-2022" moves into row 4 with updated counts (a species "Doftticka" was
# added to its species list), and the former row 4 notice "A 3379-2022" moves down into row 5 unchanged
# (other than the date update already applied above).

# New row 4: A 59219-2022
$ws.Cells.Item(4, 1).Value = "A 59219-2022"
$ws.Cells.Item(4, 2).Value = 44904
$ws.Cells.Item(4, 3).Value = 45189
$ws.Cells.Item(4, 4).Value = "VÄSTERBOTTENS LÄN"
$ws.Cells.Item(4, 5).Value = "STORUMAN"
$ws.Cells.Item(4, 6).Value = "Sveaskog"
$ws.Cells.Item(4, 7).Value = 15.6
$ws.Cells.Item(4, 8).Value = 3
$ws.Cells.Item(4, 9).Value = 3
$ws.Cells.Item(4, 10).Value = 13
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = 0
$ws.Cells.Item(4, 14).Value = 0
$ws.Cells.Item(4, 15).Value = 16
$ws.Cells.Item(4, 16).Value = 2
$ws.Cells.Item(4, 17).Value = 19
$ws.Cells.Item(4, 18).Value = "Doftticka`rFläckporing`rBlanksvart spiklav`rBlågrå svartspik`rGarnlav`rGranticka`rKolflarnlav`rKortskaftad ärgspik`rLunglav`rMörk kolflarnlav`rSkrovellav`rSpillkråka`rTalltita`rVedskivlav`rVitgrynig nållav`rMörk kådsvartspik`rBarkticka`rLuddlav`rStuplav"
$ws.Cells.Item(4, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_STORUMAN/artfynd/A 59219-2022.xlsx", "A 59219-2022")'
$ws.Cells.Item(4, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_STORUMAN/kartor/A 59219-2022.png", "A 59219-2022")'
$ws.Cells.Item(4, 21).Value = ""
$ws.Cells.Item(4, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_STORUMAN/klagomål/A 59219-2022.docx", "A 59219-2022")'
$ws.Cells.Item(4, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_STORUMAN/klagomålsmail/A 59219-2022.docx", "A 59219-2022")'
$ws.Cells.Item(4, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_STORUMAN/tillsyn/A 59219-2022.docx", "A 59219-2022")'
$ws.Cells.Item(4, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_STORUMAN/tillsynsmail/A 59219-2022.docx", "A 59219-2022")'

# New row 5: A 3379-2022 (moved down from old row 4, unchanged apart from the date already updated)
$ws.Cells.Item(5, 1).Value = "A 3379-2022"
$ws.Cells.Item(5, 2).Value = 44585
$ws.Cells.Item(5, 3).Value = 45189
$ws.Cells.Item(5, 4).Value = "VÄSTERBOTTENS LÄN"
$ws.Cells.Item(5, 5).Value = "STORUMAN"
$ws.Cells.Item(5, 6).Value = ""
$ws.Cells.Item(5, 7).Value = 3.8
$ws.Cells.Item(5, 8).Value = 7
$ws.Cells.Item(5, 9).Value = 5
$ws.Cells.Item(5, 10).Value = 11
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).Value = 0
$ws.Cells.Item(5, 14).Value = 0
$ws.Cells.Item(5, 15).Value = 11
$ws.Cells.Item(5, 16).Value = 0
$ws.Cells.Item(5, 17).Value = 18
$ws.Cells.Item(5, 18).Value = "Brunpudrad nållav`rGarnlav`rGranticka`rHarticka`rJärpe`rSkrovellav`rSpillkråka`rTalltita`rTretåig hackspett`rUllticka`rVitgrynig nållav`rBårdlav`rFinbräken`rMörkhövdad spiklav`rSpindelblomster`rStuplav`rSkogsrör`rRevlummer"
$ws.Cells.Item(5, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_STORUMAN/artfynd/A 3379-2022.xlsx", "A 3379-2022")'
$ws.Cells.Item(5, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_STORUMAN/kartor/A 3379-2022.png", "A 3379-2022")'
$ws.Cells.Item(5, 21).Value = ""
$ws.Cells.Item(5, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_STORUMAN/klagomål/A 3379-2022.docx", "A 3379-2022")'
$ws.Cells.Item(5, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_STORUMAN/klagomålsmail/A 3379-2022.docx", "A 3379-2022")'
$ws.Cells.Item(5, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_STORUMAN/tillsyn/A 3379-2022.docx", "A 3379-2022")'
$ws.Cells.Item(5, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_STORUMAN/tillsynsmail/A 3379-2022.docx", "A 3379-2022")'
